# finish admin/user/maintain edit, delete function; next: add user
#
# Sets the Default-value column (F) for several Users-table columns that
# previously had an empty default ('') and gives them real defaults:
#   email   -> 'example@email.com'
#   role    -> 'Member'
#   setting -> '{}'
#   remark  -> '{}'
# Column H recalculates automatically from the shared formula that
# references column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Leading apostrophes are Excel's "store as text" quote-prefix marker, so a
# literal leading apostrophe must be escaped by doubling it (standard Excel
# convention), otherwise .Value would silently drop it.
$ws.Range("F14").Value = "''{}'"
$ws.Range("F15").Value = "''{}'"
$ws.Range("F12").Value = "''Member'"
$ws.Range("F11").Value = "''example@email.com'"

# Move the active selection to E11 (single cell), matching the saved file.
$ws.Range("E11").Select()
